# "add real time notification" — add a "Master Content Code" header column
# to the rbt export sheet (new column J), matching it to the same header
# style as the existing G/H/I header cells, give the new column a sensible
# width, and leave the selection parked on the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1 used to hold a throwaway single-space placeholder string styled with
# a centered-but-unbordered header format; replace it with the real header
# label. Setting .Value2 rewrites the shared-string table in place (the old
# orphaned " " entry is dropped automatically once nothing references it).
$ws.Range("J1").Value2 = "Master Content Code"

# Give J1 the same visual style as the other bordered/shaded header cells
# (G1/H1/I1) instead of the old ad-hoc centered style.
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Size the new column so the longer header text isn't clipped.
$ws.Columns("J").ColumnWidth = 21

# Leave the view parked on the newly added header cell.
$ws.Range("J2").Select() | Out-Null
